$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores plain text that sometimes LOOKS like a
# simple decimal (e.g. "595.12"). Excel auto-converts such literals to a
# Number when assigned directly, so for those cells we briefly force a
# Text number format, assign the literal, then restore General -- the
# stored cell stays Text (matching the workbook's original inlineStr
# strings) without leaving a stray custom number format behind.

$ws.Range("D2").Value = "68.367.65"
$ws.Range("E2").Value = "  +1.03%  "
$ws.Range("D3").Value = "3.748.84"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.12"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.10"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("D7").Value = "3.747.38"
$ws.Range("E7").Value = "  -0.58%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  -0.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -2.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.47"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("E13").Value = "  -6.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.98"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D15").Value = "4.377.83"
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("D16").Value = "3.753.38"
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("D17").Value = "68.368.84"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.96"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -2.79%  "
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("E20").Value = "  -2.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.75"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +2.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "465.70"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.697"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -2.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.39"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("E25").Value = "  -1.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.19"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.98"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -0.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.03"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -2.37%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "3.894.71"
$ws.Range("E30").Value = "  -0.69%  "
$ws.Range("E31").Value = "  -4.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.29"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -3.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.82"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -2.12%  "
$ws.Range("E34").Value = "  -1.99%  "
$ws.Range("E35").Value = "  +0.99%  "
$ws.Range("D37").Value = "3.703.49"
$ws.Range("E37").Value = "  -0.85%  "
$ws.Range("E38").Value = "  -2.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.37"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -8.56%  "
$ws.Range("E40").Value = "  +0.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.81"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +0.44%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("E45").Value = "  -2.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "44.04"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +12.63%  "
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.91"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "145.97"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "389.58"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -1.36%  "
